$wb = $excel.ActiveWorkbook

# Sheet ALC, row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1304
$ws.Range("I19").Value = 1265.8334
$ws.Range("J19").Value = 1361.25
$ws.Range("K19").Value = 1265.8334
$ws.Range("L19").Value = 1361.25
$ws.Range("M19").Value = -1090.8334
$ws.Range("N19").Value = -1711.25

# Sheet ALC, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2671.111
$ws.Range("I40").Value = 2823
$ws.Range("J40").Value = 2367.3333
$ws.Range("K40").Value = 2823
$ws.Range("L40").Value = 2367.3333
$ws.Range("M40").Value = -2648
$ws.Range("N40").Value = -2717.3333

# Sheet ALC, row 48
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 2994.6667
$ws.Range("J48").Value = 4000
$ws.Range("L48").Value = 12000
$ws.Range("N48").Value = -12584

# Sheet ALC, row 56
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H56").Value = 2994.6667
$ws.Range("J56").Value = 4000
$ws.Range("L56").Value = 12000
$ws.Range("N56").Value = -13068

# Sheet ALC, row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 25320.334
$ws.Range("J74").Value = 41333.332
$ws.Range("L74").Value = 41333.332
$ws.Range("N74").Value = -43205.332

# Sheet ALC, row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 25320.334
$ws.Range("J77").Value = 41333.332
$ws.Range("L77").Value = 206666.66
$ws.Range("N77").Value = -216026.66

# Sheet ALC, row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2971.1667
$ws.Range("J100").Value = 4329
$ws.Range("L100").Value = 4329
$ws.Range("N100").Value = -5411

# Sheet ALC, row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 5295.75
$ws.Range("I106").Value = 2954.8
$ws.Range("K106").Value = 2954.8
$ws.Range("M106").Value = -2323.8

# Sheet ALC, row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1439.8
$ws.Range("I107").Value = 1096.5
$ws.Range("K107").Value = 1096.5
$ws.Range("M107").Value = 823.5

# Sheet ALC, row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 14601
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 14601
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 14601
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -21485

# Sheet ALC, row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 955
$ws.Range("I125").Value = 899.5
$ws.Range("J125").Value = 977.2
$ws.Range("K125").Value = 8095.5
$ws.Range("L125").Value = 8794.800000000001
$ws.Range("M125").Value = -5635.5
$ws.Range("N125").Value = -13714.8

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5653.1177
$ws.Range("I132").Value = 6212
$ws.Range("J132").Value = 3045
$ws.Range("K132").Value = 18636
$ws.Range("L132").Value = 9135
$ws.Range("M132").Value = -16106
$ws.Range("N132").Value = -14195

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2822.9565
$ws.Range("I138").Value = 2389
$ws.Range("J138").Value = 3498
$ws.Range("K138").Value = 7167
$ws.Range("L138").Value = 10494
$ws.Range("M138").Value = -2027
$ws.Range("N138").Value = -20774

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29014.29
$ws.Range("I32").Value = 29715.514
$ws.Range("K32").Value = 29715.514
$ws.Range("M32").Value = -29428.514

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1655.1428
$ws.Range("I61").Value = 1655.1428
$ws.Range("K61").Value = 1655.1428
$ws.Range("M61").Value = -1443.1428

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1655.1428
$ws.Range("I136").Value = 1655.1428
$ws.Range("K136").Value = 4965.428400000001
$ws.Range("M136").Value = -2415.428400000001

# Sheet ARM, row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 85438.336
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 85438.336
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 85438.336
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -95718.336

# Sheet BSM, row 64
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1803.3572
$ws.Range("I64").Value = 1593.5
$ws.Range("J64").Value = 2083.1667
$ws.Range("K64").Value = 1593.5
$ws.Range("L64").Value = 2083.1667
$ws.Range("M64").Value = -1368.5
$ws.Range("N64").Value = -2533.1667

# Sheet BSM, row 67
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 1803.3572
$ws.Range("I67").Value = 1593.5
$ws.Range("J67").Value = 2083.1667
$ws.Range("K67").Value = 1593.5
$ws.Range("L67").Value = 2083.1667
$ws.Range("M67").Value = -813.5
$ws.Range("N67").Value = -3643.1667

# Sheet BSM, row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4011.1304
$ws.Range("I99").Value = 3731.7896
$ws.Range("J99").Value = 5338
$ws.Range("K99").Value = 3731.7896
$ws.Range("L99").Value = 5338
$ws.Range("M99").Value = -2233.7896
$ws.Range("N99").Value = -8334

# Sheet BSM, row 130
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 74999
$ws.Range("J130").Value = 74999
$ws.Range("L130").Value = 74999
$ws.Range("N130").Value = -85039

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3767.4167
$ws.Range("I31").Value = 1991.0625
$ws.Range("J31").Value = 7320.125
$ws.Range("K31").Value = 1991.0625
$ws.Range("L31").Value = 7320.125
$ws.Range("M31").Value = -1696.0625
$ws.Range("N31").Value = -7910.125

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3767.4167
$ws.Range("I34").Value = 1991.0625
$ws.Range("J34").Value = 7320.125
$ws.Range("K34").Value = 1991.0625
$ws.Range("L34").Value = 7320.125
$ws.Range("M34").Value = -1789.0625
$ws.Range("N34").Value = -7724.125

# Sheet CRP, row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 49999
$ws.Range("J51").Value = 49999
$ws.Range("L51").Value = 49999
$ws.Range("N51").Value = -51471

# Sheet CRP, row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 49999
$ws.Range("J61").Value = 49999
$ws.Range("L61").Value = 49999
$ws.Range("N61").Value = -50695

# Sheet CRP, row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1298
$ws.Range("J94").Value = 1298
$ws.Range("L94").Value = 1298
$ws.Range("N94").Value = -2200

# Sheet CRP, row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2554.3667
$ws.Range("I107").Value = 957.53845
$ws.Range("J107").Value = 3775.4707
$ws.Range("K107").Value = 957.53845
$ws.Range("L107").Value = 3775.4707
$ws.Range("M107").Value = 962.46155
$ws.Range("N107").Value = -7615.4707

# Sheet CRP, row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2469
$ws.Range("I122").Value = 2353.875
$ws.Range("K122").Value = 7061.625
$ws.Range("M122").Value = -4611.625

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1903.4
$ws.Range("I132").Value = 1927
$ws.Range("J132").Value = 1887.6666
$ws.Range("K132").Value = 5781
$ws.Range("L132").Value = 5662.9998
$ws.Range("M132").Value = -3251
$ws.Range("N132").Value = -10722.9998

# Sheet CRP, row 140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 86616.164
$ws.Range("J140").Value = 86616.164
$ws.Range("L140").Value = 86616.164
$ws.Range("N140").Value = -96976.164

# Sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1196.6666
$ws.Range("I5").Value = 820
$ws.Range("K5").Value = 2460
$ws.Range("M5").Value = -2348

# Sheet CUL, row 33
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 306.33334
$ws.Range("I33").Value = 259.5
$ws.Range("J33").Value = 400
$ws.Range("K33").Value = 1557
$ws.Range("L33").Value = 2400
$ws.Range("M33").Value = -1274
$ws.Range("N33").Value = -2966

# Sheet CUL, row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 150
$ws.Range("J75").Value = 150
$ws.Range("L75").Value = 450
$ws.Range("N75").Value = -2446

# Sheet CUL, row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 150
$ws.Range("J78").Value = 150
$ws.Range("L78").Value = 1350
$ws.Range("N78").Value = -11334

# Sheet CUL, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1196.6666
$ws.Range("I135").Value = 820
$ws.Range("K135").Value = 7380
$ws.Range("M135").Value = -4845

# Sheet GSM, row 109
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# Sheet GSM, row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 10220
$ws.Range("I113").Value = 10163
$ws.Range("J113").Value = 10305.5
$ws.Range("K113").Value = 10163
$ws.Range("L113").Value = 10305.5
$ws.Range("M113").Value = -7993
$ws.Range("N113").Value = -14645.5

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 54797.105
$ws.Range("I132").Value = 68769.13
$ws.Range("J132").Value = 2402
$ws.Range("K132").Value = 206307.39
$ws.Range("L132").Value = 7206
$ws.Range("M132").Value = -203777.39
$ws.Range("N132").Value = -12266

# Sheet LTW, row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 57349.45
$ws.Range("I22").Value = 159571.42
$ws.Range("J22").Value = 2306.8462
$ws.Range("K22").Value = 159571.42
$ws.Range("L22").Value = 2306.8462
$ws.Range("M22").Value = -159276.42
$ws.Range("N22").Value = -2896.8462

# Sheet LTW, row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 57349.45
$ws.Range("I27").Value = 159571.42
$ws.Range("J27").Value = 2306.8462
$ws.Range("K27").Value = 159571.42
$ws.Range("L27").Value = 2306.8462
$ws.Range("M27").Value = -159464.42
$ws.Range("N27").Value = -2520.8462

# Sheet LTW, row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2264.1667
$ws.Range("I93").Value = 2428.5
$ws.Range("J93").Value = 1935.5
$ws.Range("K93").Value = 2428.5
$ws.Range("L93").Value = 1935.5
$ws.Range("M93").Value = -1180.5
$ws.Range("N93").Value = -4431.5

# Sheet LTW, row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2869.75
$ws.Range("I100").Value = 2551.7
$ws.Range("K100").Value = 2551.7
$ws.Range("M100").Value = -2010.7

# Sheet LTW, row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 87194.71000000001
$ws.Range("J133").Value = 87194.71000000001
$ws.Range("L133").Value = 87194.71000000001
$ws.Range("N133").Value = -92254.71000000001

# Sheet LTW, row 140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# Sheet LTW, row 141
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
